$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 45414161.82272212
$ws.Range("D2").Value = 9844.520545567508
$ws.Range("E2").Value = 645.3272768299601
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 45424654.90153021
